$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "37.702.64"
$ws.Range("E2").Value = "  -0.85%  "

$ws.Range("D3").Value = "2.026.20"
$ws.Range("E3").Value = "  -1.44%  "

$ws.Range("E4").Value = "  -0.12%  "

Set-TextValue "D5" "226.96"
$ws.Range("E5").Value = "  -1.44%  "

$ws.Range("E6").Value = "  -1.35%  "

Set-TextValue "D7" "59.25"
$ws.Range("E7").Value = "  +1.54%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("E9").Value = "  -0.95%  "

$ws.Range("E10").Value = "  +0.66%  "

$ws.Range("E11").Value = "  -0.35%  "

$ws.Range("E12").Value = "  -0.66%  "

$ws.Range("D13").Value = "2.326.33"
$ws.Range("E13").Value = "  -1.48%  "

Set-TextValue "D14" "20.99"
$ws.Range("E14").Value = "  +1.65%  "

$ws.Range("E15").Value = "  +0.59%  "

$ws.Range("E16").Value = "  -1.98%  "

$ws.Range("D17").Value = "2.035.60"
$ws.Range("E17").Value = "  -1.27%  "

$ws.Range("D18").Value = "37.718.90"
$ws.Range("E18").Value = "  -0.41%  "

$ws.Range("E19").Value = "  -2.25%  "

Set-TextValue "D20" "69.78"
$ws.Range("E20").Value = "  -0.11%  "

$ws.Range("E21").Value = "  -0.83%  "

Set-TextValue "D22" "224.08"
$ws.Range("E22").Value = "  -0.26%  "

Set-TextValue "D23" "0.998"
$ws.Range("E23").Value = "  -0.11%  "

Set-TextValue "D24" "2.35"
$ws.Range("E24").Value = "  -4.60%  "

$ws.Range("E25").Value = "  -1.96%  "

Set-TextValue "D26" "9.25"
$ws.Range("E26").Value = "  -0.65%  "

Set-TextValue "D27" "165.08"
$ws.Range("E27").Value = "  -0.70%  "

Set-TextValue "D28" "0.128"
$ws.Range("E28").Value = "  -3.01%  "

$ws.Range("E29").Value = "  -0.69%  "

$ws.Range("E30").Value = "  -5.06%  "

$ws.Range("E31").Value = "  +0.99%  "

$ws.Range("E32").Value = "  -2.58%  "

Set-TextValue "D33" "2.10"
$ws.Range("E33").Value = "  +5.99%  "

Set-TextValue "D34" "0.0601"
$ws.Range("E34").Value = "  -1.89%  "

Set-TextValue "D35" "4.48"
$ws.Range("E35").Value = "  -2.59%  "

Set-TextValue "D36" "6.34"
$ws.Range("E36").Value = "  +5.87%  "

$ws.Range("E37").Value = "  -4.02%  "

$ws.Range("E38").Value = "  -2.99%  "

Set-TextValue "D39" "1.00"
$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("D40").Value = "1.529.71"
$ws.Range("E40").Value = "  +3.37%  "

$ws.Range("E41").Value = "  -0.93%  "

Set-TextValue "D42" "96.77"
$ws.Range("E42").Value = "  -1.67%  "

Set-TextValue "D43" "16.71"
$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("E44").Value = "  -0.61%  "

$ws.Range("E45").Value = "  -2.70%  "

$ws.Range("E46").Value = "  +1.96%  "

$ws.Range("E47").Value = "  -1.98%  "

$ws.Range("E48").Value = "  -1.53%  "

Set-TextValue "D50" "7.09"
$ws.Range("E50").Value = "  +0.31%  "

$ws.Range("D51").Value = "2.215.76"
$ws.Range("E51").Value = "  -1.57%  "
